# Add a new "calendar" worksheet (test case for "add event" in calendar
# section) as the last tab, populate it with a title/category header row
# plus one sample row, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet ("company") so
# it lands at the end of the tab strip.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$calSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$calSheet.Name = "calendar"

# Header row.
$calSheet.Range("A1").Value = "title"
$calSheet.Range("B1").Value = "category"

# Sample data row used by the new "add event" test case.
$calSheet.Range("A2").Value = "Test title - 1"
$calSheet.Range("B2").Value = "Important"

# Highlight the header row with the same yellow fill used by every other
# sheet's header row in this workbook.
$calSheet.Range("A1:B1").Interior.Color = 65535

# Size the columns to fit their content.
$calSheet.Columns("A").ColumnWidth = 10.666666666666666
$calSheet.Columns("B").ColumnWidth = 9

# Match the sheet's recorded selection/active cell.
$calSheet.Range("C5").Select() | Out-Null

# Make the new sheet the active tab.
$calSheet.Activate()
